$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$data = @{
    2  = @(8, 8)
    3  = @(9, 9)
    4  = @(6, 7)
    5  = @(7, 7)
    6  = @(6, 7)
    7  = @(8, 8)
    8  = @(5, 6)
    9  = @(9, 9)
    10 = @(8, 8)
    11 = @(6, 7)
    12 = @(6, 6)
    13 = @(7, 7)
    14 = @(10, 10)
    15 = @(7, 7)
    16 = @(8, 8)
    17 = @(9, 9)
    18 = @(5, 5)
    19 = @(4, 5)
    20 = @(9, 9)
    21 = @(8, 8)
    22 = @(7, 7)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
